# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" positioned right after "总计" and
#    before "2021-Q3", and populate it with the quarterly fund-holding table.
# 2. Update the "总计" summary sheet with a new row for 2022-Q4 (shifting
#    the existing 2021-Q3 / 2021-Q2 rows down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" sheet before "2021-Q3"
# ---------------------------------------------------------------------
$q3Sheet   = $wb.Worksheets.Item("2021-Q3")
$newSheet  = $wb.Worksheets.Add($q3Sheet)
$newSheet.Name = "2022-Q4"

# Reference sheet used only to borrow cell formatting (style index 2:
# bold font + thin border + centered alignment) used throughout this
# workbook for header / index cells.
$styleSrc = $wb.Worksheets.Item("总计")

# -- header row (labels live in columns B:H, column A is the blank
#    pandas-index header) -------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$styleSrc.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# -- data row 2 : 012315 ---------------------------------------------------
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "012315"
$newSheet.Range("C2").Value = "创金合信港股通成长股票A"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.12"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "89.18"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "6.04"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0072"
$newSheet.Range("H2").Value = 8

# -- data row 3 : 012316 ---------------------------------------------------
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "012316"
$newSheet.Range("C3").Value = "创金合信港股通成长股票C"
$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "0.11"
$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "89.18"
$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "6.04"
$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "0.0066"
$newSheet.Range("H3").Value = 8

# index column (A2:A3) formatting
$styleSrc.Range("A2").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert the 2022-Q4 row, push
#    2021-Q3 / 2021-Q2 down one row each.
# ---------------------------------------------------------------------
$ws0 = $wb.Worksheets.Item("总计")

# Row 4 (new): 2021-Q2 data, previously on row 3
$ws0.Range("A4").Value = 2
$ws0.Range("B4").Value = "2021-Q2"
$ws0.Range("C4").Value = 3
$ws0.Range("D4").Value = 0.12
$ws0.Range("A3").Copy()
$ws0.Range("A4").PasteSpecial(-4122)

# Row 3: 2021-Q3 data, previously on row 2
$ws0.Range("A3").Value = 1
$ws0.Range("B3").Value = "2021-Q3"
$ws0.Range("C3").Value = 1
$ws0.Range("D3").Value = 0.02

# Row 2: brand-new 2022-Q4 data
$ws0.Range("A2").Value = 0
$ws0.Range("B2").Value = "2022-Q4"
$ws0.Range("C2").Value = 2
$ws0.Range("D2").Value = 0.01

# ---------------------------------------------------------------------
# Restore the originally-selected tab (2021-Q2, last sheet) as active.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
